# Apply changes described by the commit:
#   "Add files via upload / remove date contradictions."
#
# The "checks" sheet originally contained 7 contradiction-check rows:
#   row2: date_of_submission|data_provision  -> "Data provision before date of submission"      (rounddown... > 0)
#   row3: desired_dataset|available_dataset  -> "Available dataset is smaller than desired dataset" ([desired_dataset] < [available_dataset])
#   row4: date_of_submission|archiving       -> "Archiving before date of submission"            (rounddown... > 0)
#   row5: archiving|data_provision           -> "Archiving before data provision"                (rounddown... > 0)
#   row6: date_of_submission|data_provision  -> "Data provision without data submission"
#   row7: date_of_submission|archiving       -> "archiving without date of submission"
#   row8: archiving|data_provision           -> "Archiving without data provision"
#
# The three "before" (date-ordering, rounddown-based) contradiction rows
# (rows 2, 4 and 5) are removed entirely, the remaining "without" rows shift
# up, the desired/available comparison operator flips from "<" to ">", and
# the "archiving without date of submission" label is capitalised.

$wb = $excel.ActiveWorkbook
$wsItems = $wb.Worksheets.Item("item_level")
$wsChecks = $wb.Worksheets.Item("checks")

# Remove the obsolete "before" date-ordering contradiction rows.
# Deleting from the bottom up keeps the remaining row numbers stable.
$wsChecks.Rows.Item(5).Delete()
$wsChecks.Rows.Item(4).Delete()
$wsChecks.Rows.Item(2).Delete()

# Flip the comparison operator for the desired/available dataset check.
$wsChecks.Cells.Item(2, 3).Value = "[desired_dataset] > [available_dataset]"

# Capitalise the archiving-without-submission label.
$wsChecks.Cells.Item(4, 2).Value = "Archiving without date of submission"

# Update the active sheet / selections to match the saved view state:
# "item_level" is no longer the selected tab, "checks" becomes active,
# and each sheet's remembered selection moves.
$wsItems.Activate()
$wsItems.Range("C4").Select()

$wsChecks.Activate()
$wsChecks.Range("B5").Select()
